# Auto-generated edit script: updates FFXIV Tonberry Profits market-data cells
# (currentAveragePrice / LevePrice / LeveProfit columns) across all 8 sheets,
# mirroring a scheduled market-data refresh. No formulas involved; every
# touched cell is a plain numeric literal (or cleared, for 2 now-blank cells).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M51").Value = -2515
$ws.Range("H51").Value = 4309.8
$ws.Range("K51").Value = 2999
$ws.Range("I51").Value = 2999
$ws.Range("I52").Value = 1000
$ws.Range("K52").Value = 3000
$ws.Range("M52").Value = -2840
$ws.Range("J52").Value = 3988.2856
$ws.Range("L52").Value = 11964.8568
$ws.Range("H52").Value = 3324.2222
$ws.Range("N52").Value = -12284.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -1152.247
$ws.Range("H32").Value = 2015.9893
$ws.Range("K32").Value = 1439.247
$ws.Range("I32").Value = 1439.247
$ws.Range("M61").Value = -1811.9333
$ws.Range("K61").Value = 2023.9333
$ws.Range("H61").Value = 3187.56
$ws.Range("J61").Value = 4933
$ws.Range("L61").Value = 4933
$ws.Range("N61").Value = -5357
$ws.Range("I61").Value = 2023.9333
$ws.Range("K132").Value = 7715.869499999999
$ws.Range("M132").Value = -5185.869499999999
$ws.Range("N132").Value = -15620.375
$ws.Range("L132").Value = 10560.375
$ws.Range("I132").Value = 2571.9565
$ws.Range("J132").Value = 3520.125
$ws.Range("H132").Value = 2816.6453
$ws.Range("N136").Value = -19899
$ws.Range("H136").Value = 3187.56
$ws.Range("J136").Value = 4933
$ws.Range("L136").Value = 14799
$ws.Range("I136").Value = 2023.9333
$ws.Range("M136").Value = -3521.7999
$ws.Range("K136").Value = 6071.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 100
$ws.Range("H22").Value = 300
$ws.Range("N22").Value = -446
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -327
$ws.Range("L22").Value = 100
$ws.Range("I134").Value = 6030.913
$ws.Range("M134").Value = -15557.739
$ws.Range("H134").Value = 5868.44
$ws.Range("K134").Value = 18092.739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K132").Value = 5274.8181
$ws.Range("M132").Value = -2744.8181
$ws.Range("N132").Value = -22613
$ws.Range("L132").Value = 17553
$ws.Range("I132").Value = 1758.2727
$ws.Range("J132").Value = 5851
$ws.Range("H132").Value = 2635.2856
$ws.Range("J134").Value = 3783
$ws.Range("K134").Value = 3617.5716
$ws.Range("N134").Value = -16419
$ws.Range("I134").Value = 1205.8572
$ws.Range("M134").Value = -1082.5716
$ws.Range("H134").Value = 1528
$ws.Range("L134").Value = 11349
$ws.Range("H137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("L138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("N141").Value = -81020
$ws.Range("L141").Value = 70660
$ws.Range("J141").Value = 70660
$ws.Range("H141").Value = 70660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N2").Value = -598
$ws.Range("L2").Value = 372
$ws.Range("J2").Value = 62
$ws.Range("H2").Value = 317.42856
$ws.Range("H51").Value = 2146.3333
$ws.Range("N51").Value = -7670
$ws.Range("L51").Value = 6750
$ws.Range("I51").Value = 1939
$ws.Range("J51").Value = 2250
$ws.Range("K51").Value = 5817
$ws.Range("M51").Value = -5357
$ws.Range("H68").Value = 1563.4318
$ws.Range("L68").Value = 5212.9998
$ws.Range("J68").Value = 1737.6666
$ws.Range("N68").Value = -6834.9998
$ws.Range("N71").Value = -23750.9994
$ws.Range("H71").Value = 1563.4318
$ws.Range("J71").Value = 1737.6666
$ws.Range("L71").Value = 15638.9994
$ws.Range("H86").Value = 2099.6
$ws.Range("H89").Value = 2099.6
$ws.Range("J104").Value = 4027
$ws.Range("L104").Value = 12081
$ws.Range("N104").Value = -17323
$ws.Range("H104").Value = 3901.3333
$ws.Range("L131").Value = 53828.064
$ws.Range("J131").Value = 17942.688
$ws.Range("H131").Value = 9450248
$ws.Range("N131").Value = -63908.064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I122").Value = 2291.1428
$ws.Range("J122").Value = 2605.7144
$ws.Range("N122").Value = -12717.1432
$ws.Range("L122").Value = 7817.1432
$ws.Range("K122").Value = 6873.428400000001
$ws.Range("M122").Value = -4423.428400000001
$ws.Range("H122").Value = 2448.4285
$ws.Range("K132").Value = 6078662.4
$ws.Range("M132").Value = -6076132.4
$ws.Range("N132").Value = -16650.875
$ws.Range("L132").Value = 11590.875
$ws.Range("I132").Value = 2026220.8
$ws.Range("J132").Value = 3863.625
$ws.Range("H132").Value = 1427003.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L7").Value = 3832.077
$ws.Range("M7").Value = -9892
$ws.Range("J7").Value = 3832.077
$ws.Range("K7").Value = 10004
$ws.Range("I7").Value = 10004
$ws.Range("N7").Value = -4056.077
$ws.Range("H7").Value = 4272.9287
$ws.Range("M61").Value = -3532.6667
$ws.Range("K61").Value = 3734.6667
$ws.Range("H61").Value = 4180.636
$ws.Range("J61").Value = 4347.875
$ws.Range("L61").Value = 4347.875
$ws.Range("N61").Value = -4751.875
$ws.Range("I61").Value = 3734.6667
$ws.Range("H100").Value = 2691.9
$ws.Range("K100").Value = 2437.6667
$ws.Range("I100").Value = 2437.6667
$ws.Range("M100").Value = -1896.6667
$ws.Range("N113").Value = -8687.875
$ws.Range("J113").Value = 4347.875
$ws.Range("K113").Value = 3734.6667
$ws.Range("M113").Value = -1564.6667
$ws.Range("I113").Value = 3734.6667
$ws.Range("L113").Value = 4347.875
$ws.Range("H113").Value = 4180.636
$ws.Range("K126").Value = 30012
$ws.Range("I126").Value = 10004
$ws.Range("L126").Value = 11496.231
$ws.Range("N126").Value = -16436.231
$ws.Range("H126").Value = 4272.9287
$ws.Range("J126").Value = 3832.077
$ws.Range("M126").Value = -27542
$ws.Range("L132").Value = 12223.32
$ws.Range("J132").Value = 4074.44
$ws.Range("N132").Value = -17283.32
$ws.Range("H132").Value = 3780.5356

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I122").Value = 61108.78
$ws.Range("J122").Value = 1374.8334
$ws.Range("N122").Value = -9024.5002
$ws.Range("L122").Value = 4124.5002
$ws.Range("K122").Value = 183326.34
$ws.Range("M122").Value = -180876.34
$ws.Range("H122").Value = 48750.035
$ws.Range("K132").Value = 4308.6522
$ws.Range("M132").Value = -1778.6522
$ws.Range("N132").Value = -12739.4
$ws.Range("L132").Value = 7679.400000000001
$ws.Range("I132").Value = 1436.2174
$ws.Range("J132").Value = 2559.8
$ws.Range("H132").Value = 1776.697
$ws.Range("N136").Value = -15580.2633
$ws.Range("H136").Value = 11577300
$ws.Range("J136").Value = 3493.4211
$ws.Range("L136").Value = 10480.2633
$ws.Range("I136").Value = 19160138
$ws.Range("M136").Value = -57477864
$ws.Range("K136").Value = 57480414
